$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header additions for new columns H (BI_0) and I (BI_1)
$ws.Range("H1").Value = "BI_0"
$ws.Range("I1").Value = "BI_1"

# Update data rows 2-9 (A:G) with new values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 0.48352533365535499
$ws.Range("C2").Value = 0.42909300128563399
$ws.Range("D2").Value = 0.44610010709454301
$ws.Range("E2").Value = 0.97
$ws.Range("F2").Value = 0.99
$ws.Range("G2").Value = 0.97

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 0.417529865529359
$ws.Range("C3").Value = 0.40563969982942599
$ws.Range("D3").Value = 0.410637058847998
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 0.47383718367516697
$ws.Range("C4").Value = 0.44230335605386001
$ws.Range("D4").Value = 0.45061710275209998
$ws.Range("E4").Value = 0.94
$ws.Range("F4").Value = 0.99
$ws.Range("G4").Value = 0.99

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 0.500284963090742
$ws.Range("C5").Value = 0.46159578989733901
$ws.Range("D5").Value = 0.48845917918963899
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 0.36105030482545097
$ws.Range("C6").Value = 0.32340163772705199
$ws.Range("D6").Value = 0.32794953066301202
$ws.Range("E6").Value = 0.93
$ws.Range("F6").Value = 0.98
$ws.Range("G6").Value = 1

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 0.34865799944741699
$ws.Range("C7").Value = 0.333351976229248
$ws.Range("D7").Value = 0.34045231838708101
$ws.Range("E7").Value = 0.88
$ws.Range("F7").Value = 0.98
$ws.Range("G7").Value = 0.95

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 0.53755127486465404
$ws.Range("C8").Value = 0.454882532948647
$ws.Range("D8").Value = 0.49484324661663198
$ws.Range("E8").Value = 0.98
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 0.43106960812087097
$ws.Range("C9").Value = 0.40777142501032199
$ws.Range("D9").Value = 0.40681328040001002
$ws.Range("E9").Value = 0.93
$ws.Range("F9").Value = 0.99
$ws.Range("G9").Value = 0.97

# New formula columns H (BI_0 = D-B) and I (BI_1 = D-C)
$ws.Range("H2").Formula = "=D2-B2"
$ws.Range("H3:H9").Formula = "=D3-B3"
$ws.Range("I2").Formula = "=D2-C2"
$ws.Range("I3:I9").Formula = "=D3-C3"

# Update selection to match final state
$ws.Range("K9:L9").Select()
